# Revisi Gant Chart dan WBS
# - Mengganti gant chart yang lama dan WBS yang lama dengan yang baru

$wb = $excel.ActiveWorkbook

# --- Sheet "PPL" (the Gantt chart) ---
$ppl = $wb.Worksheets.Item("PPL")

$ppl.Range("B8").Value  = "Buat task,Estimasi, Software Model,  dan Sprint backlog"
$ppl.Range("B10").Value = "Update Task"
$ppl.Range("B11").Value = "Develop fitur profile"
$ppl.Range("B12").Value = "Develop fitur schedule"
$ppl.Range("B13").Value = "Develop fitur reminder"
$ppl.Range("B17").Value = "Sprint Review, Evaluation, & Update Task"
$ppl.Range("B18").Value = "Develop fitur  food calories"
$ppl.Range("B19").Value = "Develop fitur home"
$ppl.Range("B23").Value = "Sprint Review, Evaluation, & Update Task"
$ppl.Range("B24").Value = "Develop fitur week evaluation"
$ppl.Range("B25").Value = "Develop fitur food recommendation"

# --- Sheet "Discussions" ---
$disc = $wb.Worksheets.Item("Discussions")

$disc.Range("E1").Value = "27/06/13 9:02"
$disc.Range("E2").Value = "27/06/13 9:02"
